$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.012.57"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "2.356.22"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.678"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.09"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.80%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.05"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +13.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.46%  "

$ws.Range("D15").Value = "2.703.30"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.902"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").Value = "2.350.43"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").Value = "43.878.90"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +23.27%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.36%  "

$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("E40").Value = "  +4.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.202"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.24%  "

$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.91%  "

$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.10"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.74%  "
